$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1704180064308682
$ws.Range("C2").Value = 0.6205787781350482
$ws.Range("J2").Value = 0.022508038585209
$ws.Range("P2").Value = 0.1061093247588424
$ws.Range("S2").Value = 0.08038585209003216
$ws.Range("B3").Value = 0.02030456852791878
$ws.Range("C3").Value = 0.02538071065989848
$ws.Range("J3").Value = 0.03045685279187817
$ws.Range("P3").Value = 0.7868020304568528
$ws.Range("S3").Value = 0.1370558375634518
$ws.Range("P4").Value = 0.7674418604651163
$ws.Range("S4").Value = 0.2325581395348837
$ws.Range("B6").Value = 0.07207207207207207
$ws.Range("D6").Value = 0.009009009009009009
$ws.Range("F6").Value = 0.04954954954954955
$ws.Range("J6").Value = 0.2522522522522522
$ws.Range("O6").Value = 0.04054054054054054
$ws.Range("Q6").Value = 0.1756756756756757
$ws.Range("R6").Value = 0.0990990990990991
$ws.Range("S6").Value = 0.3018018018018018
$ws.Range("B7").Value = 0.115
$ws.Range("D7").Value = 0.015
$ws.Range("F7").Value = 0.075
$ws.Range("J7").Value = 0.12
$ws.Range("O7").Value = 0.015
$ws.Range("Q7").Value = 0.18
$ws.Range("R7").Value = 0.06
$ws.Range("S7").Value = 0.42
$ws.Range("B8").Value = 0.1271929824561404
$ws.Range("D8").Value = 0.01973684210526316
$ws.Range("F8").Value = 0.05263157894736842
$ws.Range("J8").Value = 0.1403508771929824
$ws.Range("O8").Value = 0.0131578947368421
$ws.Range("Q8").Value = 0.1732456140350877
$ws.Range("R8").Value = 0.1052631578947368
$ws.Range("S8").Value = 0.3684210526315789
$ws.Range("B9").Value = 0.09333333333333334
$ws.Range("D9").Value = 0.01333333333333333
$ws.Range("F9").Value = 0.07555555555555556
$ws.Range("J9").Value = 0.1111111111111111
$ws.Range("O9").Value = 0.01333333333333333
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.08888888888888889
$ws.Range("S9").Value = 0.4044444444444444
$ws.Range("B10").Value = 0.09817351598173515
$ws.Range("D10").Value = 0.02054794520547945
$ws.Range("E10").Value = 0.0015220700152207
$ws.Range("F10").Value = 0.0669710806697108
$ws.Range("J10").Value = 0.1415525114155251
$ws.Range("O10").Value = 0.0106544901065449
$ws.Range("Q10").Value = 0.2039573820395738
$ws.Range("R10").Value = 0.08904109589041095
$ws.Range("S10").Value = 0.3675799086757991
$ws.Range("G11").Value = 0.1237785016286645
$ws.Range("J11").Value = 0.07166123778501629
$ws.Range("K11").Value = 0.1693811074918567
$ws.Range("L11").Value = 0.5960912052117264
$ws.Range("S11").Value = 0.03908794788273615
$ws.Range("G12").Value = 0.6825396825396826
$ws.Range("J12").Value = 0.2222222222222222
$ws.Range("K12").Value = 0.01058201058201058
$ws.Range("L12").Value = 0.02645502645502645
$ws.Range("S12").Value = 0.0582010582010582
$ws.Range("G13").Value = 0.6862745098039216
$ws.Range("J13").Value = 0.2352941176470588
$ws.Range("S13").Value = 0.07843137254901961
$ws.Range("F15").Value = 0.02252252252252252
$ws.Range("H15").Value = 0.1576576576576577
$ws.Range("I15").Value = 0.08108108108108109
$ws.Range("J15").Value = 0.3693693693693694
$ws.Range("K15").Value = 0.05405405405405406
$ws.Range("M15").Value = 0.01801801801801802
$ws.Range("O15").Value = 0.05855855855855856
$ws.Range("S15").Value = 0.2387387387387387
$ws.Range("F16").Value = 0.02314814814814815
$ws.Range("H16").Value = 0.1759259259259259
$ws.Range("I16").Value = 0.08796296296296297
$ws.Range("J16").Value = 0.3425925925925926
$ws.Range("K16").Value = 0.09259259259259259
$ws.Range("M16").Value = 0.02777777777777778
$ws.Range("O16").Value = 0.06018518518518518
$ws.Range("S16").Value = 0.1898148148148148
$ws.Range("F17").Value = 0.01956521739130435
$ws.Range("H17").Value = 0.1869565217391304
$ws.Range("I17").Value = 0.1130434782608696
$ws.Range("J17").Value = 0.4
$ws.Range("K17").Value = 0.07391304347826087
$ws.Range("M17").Value = 0.01739130434782609
$ws.Range("N17").Value = 0.002173913043478261
$ws.Range("O17").Value = 0.0391304347826087
$ws.Range("S17").Value = 0.1478260869565217
$ws.Range("F18").Value = 0.0182648401826484
$ws.Range("H18").Value = 0.1598173515981735
$ws.Range("I18").Value = 0.1141552511415525
$ws.Range("J18").Value = 0.4063926940639269
$ws.Range("K18").Value = 0.0958904109589041
$ws.Range("M18").Value = 0.0136986301369863
$ws.Range("O18").Value = 0.0730593607305936
$ws.Range("S18").Value = 0.1187214611872146
$ws.Range("F19").Value = 0.01341281669150522
$ws.Range("H19").Value = 0.2004470938897168
$ws.Range("I19").Value = 0.08122205663189269
$ws.Range("J19").Value = 0.3450074515648286
$ws.Range("K19").Value = 0.1207153502235469
$ws.Range("M19").Value = 0.02309985096870343
$ws.Range("O19").Value = 0.07451564828614009
$ws.Range("S19").Value = 0.1415797317436662
